# Applies the crypto price/volume update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.808.82'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.22%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.910.43'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.34%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '589.39'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.43%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.50'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.85%  '

# Row 7
$ws.Range("E7").Value = '  +0.07%  '

# Row 8
$ws.Range("E8").Value = '  -0.18%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.87'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.42%  '

# Row 10
$ws.Range("E10").Value = '  -2.32%  '

# Row 11
$ws.Range("E11").Value = '  -2.29%  '

# Row 12
$ws.Range("E12").Value = '  -0.83%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '33.38'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.85%  '

# Row 14
$ws.Range("E14").Value = '  -0.24%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.392.93'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.29%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '60.740.93'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.24%  '

# Row 17
$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.67'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.68%  '

# Row 18
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.910.29'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.44%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '430.64'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.04%  '

# Row 20
$ws.Range("E20").Value = '  -2.03%  '

# Row 21
$ws.Range("E21").Value = '  -1.32%  '

# Row 22
$ws.Range("E22").Value = '  -0.72%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '81.15'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.55%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.89'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.29%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.18'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.69%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.75'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.15%  '

# Row 27
$ws.Range("E27").Value = '  +0.00%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.27'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.62%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.59'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.71%  '

# Row 30
$ws.Range("E30").Value = '  -3.66%  '

# Row 31
$ws.Range("E31").Value = '  -0.62%  '

# Row 32
$ws.Range("E32").Value = '  +1.84%  '

# Row 33
$ws.Range("E33").Value = '  +0.05%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0861'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.49%  '

# Row 35
$ws.Range("E35").Value = '  -0.41%  '

# Row 36
$ws.Range("E36").Value = '  -0.91%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.98'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.79%  '

# Row 38
$ws.Range("E38").Value = '  -1.76%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.120'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.12%  '

# Row 40
$ws.Range("E40").Value = '  -1.53%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '41.44'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.91%  '

# Row 42
$ws.Range("E42").Value = '  -5.45%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '374.53'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.25%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.694.19'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.12%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0342'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.62%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '133.35'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.43%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.77'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.38%  '

# Row 49
$ws.Range("E49").Value = '  -0.92%  '

# Row 50
$ws.Range("E50").Value = '  -3.66%  '

# Row 51
$ws.Range("E51").Value = '  -1.22%  '
